$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.6266875099344738, 1.392408230027835, -0.2280729557114796, 2.736725697326687)"
$ws.Range("C2").Value = "NIG(2.540468707914394, 1.87782560168174, 1.804583305940104, 9.286795582252218)"
$ws.Range("D2").Value = "JSU(-0.9400364873029723, 1.1432200249052609, 0.5674178242704842, 2.811983757539024)"
$ws.Range("E2").Value = "NIG(0.815473353855582, 0.5285306320037302, 5.239236279912361, 5.152398861334832)"
